$wb = $excel.ActiveWorkbook

# Update the "Last Updated" timestamp on the Metadata sheet
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "30 Oct 2025, 12:55 PM"

# Update the "Distance From Sma50" values on the "distance from Dma50" sheet
$ws = $wb.Worksheets.Item("distance from Dma50")

$ws.Range("C2").Value = 10.0454
$ws.Range("C3").Value = 7.6391
$ws.Range("C4").Value = 6.342
$ws.Range("C5").Value = 5.3792
$ws.Range("C6").Value = 5.2246
$ws.Range("C7").Value = 5.0363
$ws.Range("C8").Value = 4.4731
$ws.Range("C9").Value = 4.3545
$ws.Range("C10").Value = 3.9162
$ws.Range("C11").Value = 3.604
$ws.Range("C12").Value = 3.4175
$ws.Range("C13").Value = 3.3559
$ws.Range("C14").Value = 3.1023
$ws.Range("C15").Value = 3.0705
$ws.Range("C16").Value = 2.9878
$ws.Range("C17").Value = 2.8404
$ws.Range("C18").Value = 2.6316
$ws.Range("C19").Value = 2.5977
$ws.Range("C20").Value = 2.3863
$ws.Range("C21").Value = 2.3494
$ws.Range("C22").Value = 1.4418
$ws.Range("C23").Value = 1.4379
$ws.Range("C24").Value = 1.3802
$ws.Range("C25").Value = 1.2239
$ws.Range("C26").Value = 1.0493
$ws.Range("C27").Value = 0.9487
$ws.Range("C28").Value = 0.6157
$ws.Range("C29").Value = 0.3083
$ws.Range("C30").Value = -2.0744
